# Update order status/payment on the "All Orders" sheet
$wb = $excel.ActiveWorkbook

$wsOrders = $wb.Worksheets.Item("All Orders")
$wsOrders.Range("H2").Value = "DELIVERED"
$wsOrders.Range("I2").Value = "PAID"

# Update the corresponding daily rollup on the "Daily Summary" sheet
$wsSummary = $wb.Worksheets.Item("Daily Summary")
$wsSummary.Range("C2").Value = 1
$wsSummary.Range("F2").Value = 120
$wsSummary.Range("G2").Value = 405
